$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.6060606060606061
$ws.Range("F2").Value = 0.6149903569258408
$ws.Range("G2").Value = 0.6060606060606061
$ws.Range("H2").Value = 0.5996968755589445
$ws.Range("E3").Value = 0.5959595959595959
$ws.Range("F3").Value = 0.6054994388327721
$ws.Range("G3").Value = 0.5959595959595959
$ws.Range("H3").Value = 0.5883251115809256
$ws.Range("E4").Value = 0.5757575757575758
$ws.Range("F4").Value = 0.5760546642899584
$ws.Range("G4").Value = 0.5757575757575758
$ws.Range("H4").Value = 0.5756709956709957
$ws.Range("E5").Value = 0.5656565656565656
$ws.Range("F5").Value = 0.5700725200725201
$ws.Range("G5").Value = 0.5656565656565656
$ws.Range("H5").Value = 0.5607315292914716
$ws.Range("E6").Value = 0.5656565656565656
$ws.Range("F6").Value = 0.5667789001122335
$ws.Range("G6").Value = 0.5656565656565656
$ws.Range("H6").Value = 0.5647687916466659
$ws.Range("E7").Value = 0.5858585858585859
$ws.Range("F7").Value = 0.5958694083694084
$ws.Range("G7").Value = 0.5858585858585859
$ws.Range("H7").Value = 0.5768039583760562
$ws.Range("E8").Value = 0.6464646464646465
$ws.Range("F8").Value = 0.6380163187943508
$ws.Range("G8").Value = 0.6464646464646465
$ws.Range("H8").Value = 0.6210250654695099
$ws.Range("E9").Value = 0.6767676767676768
$ws.Range("F9").Value = 0.6723717682621791
$ws.Range("G9").Value = 0.6767676767676768
$ws.Range("H9").Value = 0.6596265687174777
$ws.Range("E10").Value = 0.5858585858585859
$ws.Range("F10").Value = 0.5644007644007645
$ws.Range("G10").Value = 0.5858585858585859
$ws.Range("H10").Value = 0.5614268772163509
$ws.Range("E13").Value = 0.6060606060606061
$ws.Range("F13").Value = 0.5845615408025738
$ws.Range("G13").Value = 0.6060606060606061
$ws.Range("H13").Value = 0.5587114142545274
$ws.Range("E14").Value = 0.6363636363636364
$ws.Range("F14").Value = 0.6361399060332612
$ws.Range("G14").Value = 0.6363636363636364
$ws.Range("H14").Value = 0.6359916679065615
$ws.Range("E15").Value = 0.5656565656565656
$ws.Range("F15").Value = 0.5659492888064316
$ws.Range("G15").Value = 0.5656565656565656
$ws.Range("H15").Value = 0.5657452343172461
$ws.Range("F16").Value = 0.6705128205128206
$ws.Range("H16").Value = 0.6629388008698355
$ws.Range("E17").Value = 0.6363636363636364
$ws.Range("F17").Value = 0.6727272727272726
$ws.Range("G17").Value = 0.6363636363636364
$ws.Range("H17").Value = 0.6103896103896104
$ws.Range("E18").Value = 0.5757575757575758
$ws.Range("F18").Value = 0.5753018910913648
$ws.Range("G18").Value = 0.5757575757575758
$ws.Range("H18").Value = 0.5734006734006734
$ws.Range("E19").Value = 0.6161616161616161
$ws.Range("F19").Value = 0.6173600410888547
$ws.Range("G19").Value = 0.6161616161616161
$ws.Range("H19").Value = 0.61267217630854
$ws.Range("E20").Value = 0.7373737373737373
$ws.Range("F20").Value = 0.6736596736596737
$ws.Range("G20").Value = 0.7373737373737373
$ws.Range("H20").Value = 0.6843738590726542
$ws.Range("E21").Value = 0.8080808080808081
$ws.Range("F21").Value = 0.7943350289117864
$ws.Range("G21").Value = 0.8080808080808081
$ws.Range("H21").Value = 0.7959759157522091
$ws.Range("E22").Value = 0.6464646464646465
$ws.Range("F22").Value = 0.6063432024694578
$ws.Range("G22").Value = 0.6464646464646465
$ws.Range("H22").Value = 0.6241661605961747
$ws.Range("E24").Value = 0.7373737373737373
$ws.Range("F24").Value = 0.6931382793451759
$ws.Range("G24").Value = 0.7373737373737373
$ws.Range("H24").Value = 0.7033295922184811
$ws.Range("E25").Value = 0.7474747474747475
$ws.Range("F25").Value = 0.6772404900064475
$ws.Range("G25").Value = 0.7474747474747475
$ws.Range("H25").Value = 0.678946164357305
$ws.Range("E26").Value = 0.6262626262626263
$ws.Range("F26").Value = 0.6675084175084175
$ws.Range("G26").Value = 0.6262626262626263
$ws.Range("H26").Value = 0.5996402379964023
$ws.Range("E27").Value = 0.6060606060606061
$ws.Range("F27").Value = 0.6106622741652941
$ws.Range("G27").Value = 0.6060606060606061
$ws.Range("H27").Value = 0.6006759110207387
$ws.Range("E28").Value = 0.6363636363636364
$ws.Range("F28").Value = 0.6394219741570457
$ws.Range("G28").Value = 0.6363636363636364
$ws.Range("H28").Value = 0.6336700336700336
$ws.Range("E31").Value = 0.5959595959595959
$ws.Range("F31").Value = 0.6384356384356384
$ws.Range("G31").Value = 0.5959595959595959
$ws.Range("H31").Value = 0.5593434343434344
